$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H31").Value = 2334.1667
$ws.Range("I31").Value = 1546.4546
$ws.Range("J31").Value = 10999
$ws.Range("K31").Value = 4639.3638
$ws.Range("L31").Value = 32997
$ws.Range("M31").Value = -4409.3638
$ws.Range("N31").Value = -33457
$ws.Range("H43").Value = 12663.286
$ws.Range("J43").Value = 9548.333000000001
$ws.Range("L43").Value = 9548.333000000001
$ws.Range("N43").Value = -9686.333000000001
$ws.Range("H98").Value = 1007.2222
$ws.Range("I98").Value = 976.875
$ws.Range("J98").Value = 1250
$ws.Range("K98").Value = 976.875
$ws.Range("L98").Value = 1250
$ws.Range("M98").Value = 521.125
$ws.Range("N98").Value = -4246
$ws.Range("H105").Value = 69999.5
$ws.Range("J105").Value = 69999.5
$ws.Range("L105").Value = 69999.5
$ws.Range("N105").Value = -76987.5
$ws.Range("H106").Value = 1819.75
$ws.Range("I106").Value = 1819.75
$ws.Range("K106").Value = 1819.75
$ws.Range("M106").Value = -1188.75
$ws.Range("H107").Value = 2056.7
$ws.Range("I107").Value = 1951.8889
$ws.Range("K107").Value = 1951.8889
$ws.Range("M107").Value = -31.88889999999992
$ws.Range("H115").Value = 1590.3334
$ws.Range("I115").Value = 1590.3334
$ws.Range("J115").Value = 0
$ws.Range("K115").Value = 4771.0002
$ws.Range("L115").Value = 0
$ws.Range("N115").Value = -3204.0002
$ws.Range("H122").Value = 1007.2222
$ws.Range("I122").Value = 976.875
$ws.Range("J122").Value = 1250
$ws.Range("K122").Value = 2930.625
$ws.Range("L122").Value = 3750
$ws.Range("M122").Value = -480.625
$ws.Range("N122").Value = -8650
$ws.Range("H131").Value = 1700
$ws.Range("I131").Value = 395
$ws.Range("K131").Value = 1185
$ws.Range("M131").Value = 3855
$ws.Range("H132").Value = 3563.4102
$ws.Range("I132").Value = 3635.457
$ws.Range("J132").Value = 2933
$ws.Range("K132").Value = 10906.371
$ws.Range("L132").Value = 8799
$ws.Range("M132").Value = -8376.370999999999
$ws.Range("N132").Value = -13859
$ws.Range("H137").Value = 1781.8
$ws.Range("I137").Value = 930.5
$ws.Range("J137").Value = 2349.3333
$ws.Range("K137").Value = 2791.5
$ws.Range("L137").Value = 7047.999899999999
$ws.Range("M137").Value = -241.5
$ws.Range("N137").Value = -12147.9999
$ws.Range("H138").Value = 3293.327
$ws.Range("I138").Value = 1850.7916
$ws.Range("J138").Value = 4529.7856
$ws.Range("K138").Value = 5552.3748
$ws.Range("L138").Value = 13589.3568
$ws.Range("M138").Value = -412.3747999999996
$ws.Range("N138").Value = -23869.3568
$ws.Range("M115").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 3569.1304
$ws.Range("I132").Value = 3569.1304
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 10707.3912
$ws.Range("L132").Value = 0
$ws.Range("N132").Value = -8177.3912
$ws.Range("M132").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H106").Value = 49999
$ws.Range("J106").Value = 49999
$ws.Range("L106").Value = 49999
$ws.Range("N106").Value = -52523
$ws.Range("H133").Value = 16999.7
$ws.Range("J133").Value = 16999.7
$ws.Range("L133").Value = 16999.7
$ws.Range("N133").Value = -27119.7

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H9").Value = 49999
$ws.Range("J9").Value = 49999
$ws.Range("L9").Value = 49999
$ws.Range("N9").Value = -50335
$ws.Range("H31").Value = 7468.696
$ws.Range("I31").Value = 8601.684999999999
$ws.Range("J31").Value = 2087
$ws.Range("K31").Value = 8601.684999999999
$ws.Range("L31").Value = 2087
$ws.Range("M31").Value = -8306.684999999999
$ws.Range("N31").Value = -2677
$ws.Range("H34").Value = 7468.696
$ws.Range("I34").Value = 8601.684999999999
$ws.Range("J34").Value = 2087
$ws.Range("K34").Value = 8601.684999999999
$ws.Range("L34").Value = 2087
$ws.Range("M34").Value = -8399.684999999999
$ws.Range("N34").Value = -2491
$ws.Range("H43").Value = 49547.8
$ws.Range("J43").Value = 49547.8
$ws.Range("L43").Value = 49547.8
$ws.Range("N43").Value = -49915.8
$ws.Range("H75").Value = 48798
$ws.Range("J75").Value = 57996.668
$ws.Range("L75").Value = 57996.668
$ws.Range("N75").Value = -59992.668
$ws.Range("H78").Value = 48798
$ws.Range("J78").Value = 57996.668
$ws.Range("L78").Value = 173990.004
$ws.Range("N78").Value = -183974.004
$ws.Range("H101").Value = 49547.8
$ws.Range("J101").Value = 49547.8
$ws.Range("L101").Value = 49547.8
$ws.Range("N101").Value = -56037.8
$ws.Range("H105").Value = 1692.25
$ws.Range("I105").Value = 1692.25
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 1692.25
$ws.Range("L105").Value = 0
$ws.Range("N105").Value = 54.75
$ws.Range("H132").Value = 4590.483
$ws.Range("I132").Value = 3730.28
$ws.Range("K132").Value = 11190.84
$ws.Range("M132").Value = -8660.84
$ws.Range("H133").Value = 80073
$ws.Range("J133").Value = 84998.664
$ws.Range("L133").Value = 84998.664
$ws.Range("N133").Value = -90058.664
$ws.Range("H134").Value = 3473.853
$ws.Range("I134").Value = 3446.7812
$ws.Range("K134").Value = 10340.3436
$ws.Range("M134").Value = -7805.3436
$ws.Range("H137").Value = 62477.2
$ws.Range("J137").Value = 62221.332
$ws.Range("L137").Value = 62221.332
$ws.Range("N137").Value = -72421.33199999999
$ws.Range("M105").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H35").Value = 1000
$ws.Range("J35").Value = 1000
$ws.Range("L35").Value = 3000
$ws.Range("N35").Value = -3576
$ws.Range("H37").Value = 90970390
$ws.Range("J37").Value = 90970390
$ws.Range("L37").Value = 272911170
$ws.Range("N37").Value = -272911394
$ws.Range("H50").Value = 620.8461
$ws.Range("I50").Value = 510.6
$ws.Range("J50").Value = 988.3333
$ws.Range("K50").Value = 1531.8
$ws.Range("L50").Value = 2964.9999
$ws.Range("M50").Value = -1050.8
$ws.Range("N50").Value = -3926.9999
$ws.Range("H53").Value = 620.8461
$ws.Range("I53").Value = 510.6
$ws.Range("J53").Value = 988.3333
$ws.Range("K53").Value = 1531.8
$ws.Range("L53").Value = 2964.9999
$ws.Range("M53").Value = -1050.8
$ws.Range("N53").Value = -3926.9999
$ws.Range("H64").Value = 6130.3335
$ws.Range("I64").Value = 6255
$ws.Range("J64").Value = 6005.6665
$ws.Range("K64").Value = 18765
$ws.Range("L64").Value = 18016.9995
$ws.Range("M64").Value = -18495
$ws.Range("N64").Value = -18556.9995
$ws.Range("H67").Value = 6130.3335
$ws.Range("I67").Value = 6255
$ws.Range("J67").Value = 6005.6665
$ws.Range("K67").Value = 18765
$ws.Range("L67").Value = 18016.9995
$ws.Range("M67").Value = -17829
$ws.Range("N67").Value = -19888.9995
$ws.Range("H109").Value = 2654.9
$ws.Range("I109").Value = 1721.2858
$ws.Range("K109").Value = 5163.857400000001
$ws.Range("M109").Value = -4123.857400000001
$ws.Range("H114").Value = 945.6
$ws.Range("J114").Value = 0
$ws.Range("L114").Value = 0
$ws.Range("H115").Value = 799.5
$ws.Range("I115").Value = 799.5
$ws.Range("K115").Value = 2398.5
$ws.Range("M115").Value = -1223.5
$ws.Range("H131").Value = 2057.926
$ws.Range("I131").Value = 1392.2858
$ws.Range("J131").Value = 2290.9
$ws.Range("K131").Value = 4176.857400000001
$ws.Range("L131").Value = 6872.700000000001
$ws.Range("M131").Value = 863.1425999999992
$ws.Range("N131").Value = -16952.7
$ws.Range("H132").Value = 0
$ws.Range("I132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("N114").ClearContents()
$ws.Range("M132").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6000
$ws.Range("I70").Value = 6000
$ws.Range("K70").Value = 6000
$ws.Range("M70").Value = -5730
$ws.Range("H73").Value = 6000
$ws.Range("I73").Value = 6000
$ws.Range("K73").Value = 6000
$ws.Range("M73").Value = -5064
$ws.Range("H105").Value = 48597.4
$ws.Range("J105").Value = 48597.4
$ws.Range("L105").Value = 48597.4
$ws.Range("N105").Value = -55585.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 3737.25
$ws.Range("I16").Value = 2974.5
$ws.Range("J16").Value = 4500
$ws.Range("K16").Value = 2974.5
$ws.Range("L16").Value = 4500
$ws.Range("M16").Value = -2804.5
$ws.Range("N16").Value = -4840
$ws.Range("H22").Value = 2427.6667
$ws.Range("I22").Value = 3383.5
$ws.Range("J22").Value = 1949.75
$ws.Range("K22").Value = 3383.5
$ws.Range("L22").Value = 1949.75
$ws.Range("M22").Value = -3088.5
$ws.Range("N22").Value = -2539.75
$ws.Range("H27").Value = 2427.6667
$ws.Range("I27").Value = 3383.5
$ws.Range("J27").Value = 1949.75
$ws.Range("K27").Value = 3383.5
$ws.Range("L27").Value = 1949.75
$ws.Range("M27").Value = -3276.5
$ws.Range("N27").Value = -2163.75
$ws.Range("H55").Value = 266
$ws.Range("I55").Value = 240.25
$ws.Range("J55").Value = 317.5
$ws.Range("K55").Value = 240.25
$ws.Range("L55").Value = 317.5
$ws.Range("M55").Value = -67.25
$ws.Range("N55").Value = -663.5
$ws.Range("H127").Value = 69999
$ws.Range("J127").Value = 69999
$ws.Range("L127").Value = 69999
$ws.Range("N127").Value = -79919
$ws.Range("H138").Value = 99992.5
$ws.Range("J138").Value = 99992.5
$ws.Range("L138").Value = 99992.5
$ws.Range("N138").Value = -110272.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3511.889
$ws.Range("I132").Value = 3593.2856
$ws.Range("K132").Value = 10779.8568
$ws.Range("M132").Value = -8249.856800000001
